# Fix up a couple classes in Met Council results: rename class "XSB" to "XB"
# for the three affected rows (21-23), and update the selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename class abbreviation "XSB" -> "XB" for rows 21-23 in column B
$ws.Range("B21").Value = "XB"
$ws.Range("B22").Value = "XB"
$ws.Range("B23").Value = "XB"

# Update the selection state to match the new active selection
$ws.Range("B22:B23").Select()

$wb.Save()
